# Adding THree Evalustion metrics for Tuning
# Insert two new metric columns (calinski_harabasz_score, davies_bouldin_score)
# between silhouette_score and parameters_combinations, re-sort the rows by
# file (test_1..test_4) and by linkage (ward, complete, average, complete,
# average) to match the new evaluation ordering, and fill in all values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("C1").Value = "calinski_harabasz_score"
$ws.Range("D1").Value = "davies_bouldin_score"
$ws.Range("E1").Value = "parameters_combinations"

# Copy the header style (bold, bordered, centered) from an existing header
# cell onto the newly used E1 header cell.
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "parameters_combinations"

# --- Data rows ----------------------------------------------------------
# columns: A=file  B=silhouette_score  C=calinski_harabasz_score
#          D=davies_bouldin_score      E=parameters_combinations
$rows = @(
    @{ Row=2;  File="eval_datasets/test_1_eval.xlsx"; B=0.2789490627547889; C=7821.728977223685; D=1.06287242850952;  E='"3,euclidean,ward"' },
    @{ Row=3;  File="eval_datasets/test_1_eval.xlsx"; B=0.4605969901439101; C=2200.082270199993; D=0.6030689835627047; E='"3,euclidean,complete"' },
    @{ Row=4;  File="eval_datasets/test_1_eval.xlsx"; B=0.7900151687435545; C=539.5168011753614;  D=0.2998735945137385; E='"3,euclidean,average"' },
    @{ Row=5;  File="eval_datasets/test_1_eval.xlsx"; B=0.7524777728335076; C=660.7881165843683;  D=0.4057289160622577; E='"3,manhattan,complete"' },
    @{ Row=6;  File="eval_datasets/test_1_eval.xlsx"; B=0.7450514926598724; C=610.5149829166194;  D=0.3823409017911434; E='"3,manhattan,average"' },
    @{ Row=7;  File="eval_datasets/test_2_eval.xlsx"; B=0.3297529116811677; C=6120.725641662055;  D=0.9671650765148169; E='"3,euclidean,ward"' },
    @{ Row=8;  File="eval_datasets/test_2_eval.xlsx"; B=0.5599471565244205; C=1206.634924035146;  D=0.5206835087392979; E='"3,euclidean,complete"' },
    @{ Row=9;  File="eval_datasets/test_2_eval.xlsx"; B=0.7336845697308644; C=535.536293877974;   D=0.3691751804849007; E='"3,euclidean,average"' },
    @{ Row=10; File="eval_datasets/test_2_eval.xlsx"; B=0.5649826201906503; C=1188.901943034108;  D=0.5611769956664775; E='"3,manhattan,complete"' },
    @{ Row=11; File="eval_datasets/test_2_eval.xlsx"; B=0.7873551073615996; C=287.0972113100034;  D=0.305926278826105;  E='"3,manhattan,average"' },
    @{ Row=12; File="eval_datasets/test_3_eval.xlsx"; B=0.3385353445223021; C=8223.908197639621;  D=0.8923090917525595; E='"3,euclidean,ward"' },
    @{ Row=13; File="eval_datasets/test_3_eval.xlsx"; B=0.3539922627492189; C=1712.547115344986;  D=0.6302781897366564; E='"3,euclidean,complete"' },
    @{ Row=14; File="eval_datasets/test_3_eval.xlsx"; B=0.6790248925233603; C=199.5744768548442;  D=0.3691286926736267; E='"3,euclidean,average"' },
    @{ Row=15; File="eval_datasets/test_3_eval.xlsx"; B=0.4842450027645873; C=1323.673453672702;  D=0.5203162854306486; E='"3,manhattan,complete"' },
    @{ Row=16; File="eval_datasets/test_3_eval.xlsx"; B=0.6730601494851508; C=206.5566615213655;  D=0.3994247883164554; E='"3,manhattan,average"' },
    @{ Row=17; File="eval_datasets/test_4_eval.xlsx"; B=0.353192212347119;  C=6366.302243950843;  D=0.8081520445489757; E='"3,euclidean,ward"' },
    @{ Row=18; File="eval_datasets/test_4_eval.xlsx"; B=0.5493149780355437; C=2540.581358420732;  D=0.5246664979873296; E='"3,euclidean,complete"' },
    @{ Row=19; File="eval_datasets/test_4_eval.xlsx"; B=0.6090192210379256; C=1680.006892396578;  D=0.4463743192296767; E='"3,euclidean,average"' },
    @{ Row=20; File="eval_datasets/test_4_eval.xlsx"; B=0.362795979641472;  C=5775.863201932922;  D=0.757796438139723;  E='"3,manhattan,complete"' },
    @{ Row=21; File="eval_datasets/test_4_eval.xlsx"; B=0.6637254014072674; C=1972.434548032037;  D=0.4632990372470223; E='"3,manhattan,average"' }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.File
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
}

Write-Host "Updated range now $($ws.UsedRange.Address())"
